$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp caption
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 02:17"

# Full refreshed & re-sorted (descending by Casos totales) country data, rows 4-219
$data = @(
    @('Estados Unidos',6334593,43856,3572719,2570847,0,1063,191027),
    @('Brasil',4046150,44728,3247610,673811,0,830,124729),
    @('India',3933124,84156,3032916,831639,0,1083,68569),
    @('Rusia',1009995,4995,826935,165532,0,114,17528),
    @('Peru',670145,6708,489886,150854,0,146,29405),
    @('Colombia',641574,8235,489151,131805,0,270,20618),
    @('Sudafrica',633015,2420,554887,63565,0,174,14563),
    @('Mexico',610957,4921,424990,120151,0,575,65816),
    @('España',488513,8959,0,0,0,40,29234),
    @('Argentina',451198,12026,322461,119376,0,243,9361),
    @('Chile',416501,1762,389409,15670,0,78,11422),
    @('Iran',380746,1994,328595,30225,0,129,21926),
    @('Reino Unido',340411,1735,0,0,0,13,41527),
    @('Banglades',319686,2158,213980,101323,0,32,4383),
    @('Arabia Saudita',318319,833,293964,20373,0,26,3982),
    @('Francia',300181,7157,87206,182269,0,20,30706),
    @('Pakistan',297014,424,281925,8761,0,10,6328),
    @('Turquia',274943,1642,248087,20345,0,49,6511),
    @('Italia',272912,1397,208490,28915,0,10,35507),
    @('Alemania',248814,1423,223100,16315,0,6,9399),
    @('Irak',247039,4755,187757,52007,0,74,7275),
    @('Filipinas',228403,1987,159475,65240,0,65,3688),
    @('Indonesia',184268,3622,132055,44463,0,134,7750),
    @('Canada',130493,570,115444,5908,0,6,9141),
    @('Ucrania',128228,2430,59676,65842,0,54,2710),
    @('Israel',124455,2991,98645,24825,0,16,985),
    @('Catar',119420,214,116313,2906,0,2,201),
    @('Bolivia',117928,661,64074,48651,0,102,5203),
    @('Ecuador',116360,903,102051,7661,0,29,6648),
    @('Kazajistan',106032,88,97967,6477,0,0,1588),
    @('Egipto',99425,145,75415,18531,0,18,5479),
    @('Republica Dominicana',96629,1002,70871,23957,0,36,1801),
    @('Panama',94914,830,68198,24670,0,16,2046),
    @('Rumania',91256,1365,39275,48216,0,44,3765),
    @('Kuwait',87378,900,78791,8051,0,1,536),
    @('Oman',86380,256,81828,3847,0,16,705),
    @('Belgica',85911,424,18490,57523,0,1,9898),
    @('China',85077,11,80251,192,0,0,4634),
    @('Suecia',84729,0,0,0,0,5,5832),
    @('Guatemala',76358,714,64399,9155,0,14,2804),
    @('Paises Bajos',72464,601,0,0,0,0,6235),
    @('Bielorrusia',72302,161,71205,401,0,5,696),
    @('Emiratos Arabes Unidos',72154,614,62668,9099,0,0,387),
    @('Japon',69599,598,59524,8756,0,12,1319),
    @('Polonia',69129,612,48593,18444,0,14,2092),
    @('Marruecos',66855,1402,51223,14379,0,37,1253),
    @('Honduras',62526,757,11147,49455,0,36,1924),
    @('Portugal',59051,418,42427,14795,0,2,1829),
    @('Singapur',56908,48,56028,853,0,0,27),
    @('Etiopia',55213,804,20283,34074,0,10,856),
    @('Nigeria',54588,125,42627,10913,0,21,1048),
    @('Barein',53433,626,50013,3230,0,0,190),
    @('Venezuela',48883,0,39912,8573,0,0,398),
    @('Argelia',45469,311,32006,11934,0,6,1529),
    @('Ghana',44713,55,43577,856,0,4,280),
    @('Costa Rica',44458,1153,17855,26143,0,7,460),
    @('Armenia',44271,196,38855,4529,0,3,887),
    @('Kirguistan',44135,99,39174,3902,0,0,1059),
    @('Suiza',43127,364,36500,4614,0,2,2013),
    @('Nepal',42877,1228,24207,18413,0,6,257),
    @('Uzbekistan',42688,251,40081,2276,0,4,331),
    @('Moldavia',38372,632,27017,10319,0,12,1036),
    @('Afganistan',38288,45,29390,7489,0,0,1409),
    @('Azerbaiyan',36899,167,34116,2242,0,3,541),
    @('Kenia',34705,212,20644,13476,0,4,585),
    @('Serbia',31676,95,30188,770,0,2,718),
    @('Irlanda',29206,92,23364,4065,0,0,1777),
    @('Austria',28372,403,24171,3466,0,1,735),
    @('Chequia',26452,679,18663,7363,0,1,426),
    @('Australia',26049,126,21912,3459,0,15,678),
    @('El Salvador',26000,96,15119,10142,0,8,739),
    @('Estado de Palestina',24471,596,16095,8209,0,5,167),
    @('Bosnia y Herzegovina',20804,287,14120,6048,0,9,636),
    @('Corea del Sur',20644,195,15529,4786,0,3,329),
    @('Camerun',19604,144,18448,741,0,0,415),
    @('Paraguay',19138,0,9981,8799,0,0,358),
    @('Libano',18963,588,5338,13446,0,2,179),
    @('Costa de Marfil',18208,47,17045,1044,0,2,119),
    @('Dinamarca',17374,179,15499,1249,0,0,626),
    @('Bulgaria',16775,158,11935,4182,0,10,658),
    @('Libia',15773,617,1856,13663,0,4,254),
    @('Madagascar',15106,83,14031,878,0,1,197),
    @('Republica de Macedonia',14762,162,11965,2191,0,0,606),
    @('Senegal',13826,83,9553,3986,0,0,287),
    @('Sudan',13189,0,6612,5754,0,0,823),
    @('Zambia',12523,108,11562,669,0,0,292),
    @('Noruega',11120,86,9348,1508,0,0,264),
    @('Croacia',11094,369,8266,2634,0,3,194),
    @('Grecia',10998,241,3804,6916,0,5,278),
    @('Consejo Danes para los Refugiados',10125,11,9367,499,0,0,259),
    @('Albania',9844,116,5732,3811,0,5,301),
    @('Guinea',9579,53,8726,792,0,1,61),
    @('Malasia',9374,14,9083,163,0,0,128),
    @('Guayana Francesa',9251,42,8772,418,0,0,61),
    @('Tayikistan',8690,36,7482,1139,0,0,69),
    @('Gabon',8538,0,7335,1150,0,0,53),
    @('Haiti',8301,43,5870,2221,0,4,210),
    @('Maldivas',8281,141,5483,2769,0,0,29),
    @('Finlandia',8200,39,7350,514,0,0,336),
    @('Namibia',8082,238,3483,4513,0,4,86),
    @('Mauritania',7106,17,6588,358,0,0,160),
    @('Hungria',6923,301,3930,2373,0,1,620),
    @('Luxemburgo',6811,66,0,0,0,0,124),
    @('Zimbabue',6678,40,5263,1209,0,0,206),
    @('Malaui',5593,14,3516,1902,0,0,175),
    @('Republica de Yibuti',5387,0,5327,0,0,0,60),
    @('Montenegro',5165,146,4164,897,0,0,104),
    @('Guinea Ecuatorial',4965,0,4390,492,0,0,83),
    @('Hong Kong',4839,8,4431,314,0,2,94),
    @('Republica de Africa Central',4729,17,1807,2860,0,0,62),
    @('Suazilandia',4720,52,3789,837,0,0,94),
    @('Nicaragua',4668,0,2913,1614,0,0,141),
    @('Congo',4628,0,1742,2784,0,0,102),
    @('Tunez',4394,198,1681,2629,0,3,84),
    @('Ruanda',4255,37,2163,2074,0,1,18),
    @('Surinam',4215,66,3318,824,0,1,73),
    @('Cuba',4214,88,3474,640,0,2,100),
    @('Mozambique',4207,90,2370,1811,0,1,26),
    @('Eslovaquia',4163,121,2617,1509,0,4,37),
    @('Cabo Verde',4125,77,3502,582,0,0,41),
    @('Tailandia',3427,2,3277,92,0,0,58),
    @('Mayotte',3374,73,2964,370,0,0,40),
    @('Somalia',3310,0,2481,731,0,0,98),
    @('Uganda',3288,176,1532,1723,0,1,33),
    @('Sri Lanka',3111,10,2883,216,0,0,12),
    @('Gambia',3101,34,1075,1927,0,2,99),
    @('Eslovenia',3032,53,2402,496,0,0,134),
    @('Lituania',2978,20,1901,991,0,0,86),
    @('Siria',2973,75,681,2168,0,4,124),
    @('Jamaica',2822,139,900,1895,0,3,27),
    @('Mali',2807,5,2203,478,0,0,126),
    @('Angola',2805,28,1144,1548,0,1,113),
    @('Sudan del Sur',2533,1,1290,1196,0,0,47),
    @('Estonia',2441,26,2150,227,0,0,64),
    @('Bahamas',2337,0,849,1438,0,0,50),
    @('Aruba',2292,81,1031,1248,0,1,13),
    @('Jordania',2233,72,1648,570,0,0,15),
    @('Guinea-Bisau',2205,0,1127,1044,0,0,34),
    @('Benin',2194,0,1738,416,0,0,40),
    @('Islandia',2128,7,2023,95,0,0,10),
    @('Sierra Leona',2035,6,1600,364,0,0,71),
    @('Trinidad yTobago',1984,64,700,1255,0,1,29),
    @('Yemen',1979,3,1180,228,0,0,571),
    @('Malta',1965,34,1528,424,0,0,13),
    @('Reunion',1912,116,880,1022,0,0,10),
    @('Nueva Zelanda',1759,2,1622,115,0,0,22),
    @('Botsuana',1724,0,493,1225,0,0,6),
    @('Uruguay',1636,10,1437,155,0,0,44),
    @('Georgia',1568,20,1279,270,0,0,19),
    @('Republica de Chipre',1498,3,1139,338,0,0,21),
    @('Togo',1443,9,1071,341,0,1,31),
    @('Letonia',1410,4,1187,189,0,0,34),
    @('Guyana',1401,19,805,552,0,3,44),
    @('Burkina Faso',1375,5,1075,245,0,0,55),
    @('Guadalupe',1363,0,336,1009,0,0,18),
    @('Liberia',1306,1,1163,61,0,0,82),
    @('Principado de Andorra',1199,0,909,237,0,0,53),
    @('Niger',1177,1,1091,17,0,0,69),
    @('Birmania',1111,116,359,746,0,0,6),
    @('Belice',1101,0,240,848,0,0,13),
    @('Lesoto',1085,0,528,526,0,0,31),
    @('Vietnam',1046,0,755,256,0,1,35),
    @('Republica del Chad',1018,1,908,33,0,0,77),
    @('Santo Tome y Principe',897,1,857,25,0,0,15),
    @('Martinica',754,0,98,638,0,0,18),
    @('San Marino',735,0,660,33,0,0,42),
    @('Crucero',712,0,651,48,0,0,13),
    @('Polinesia Francesa',622,0,336,286,0,0,0),
    @('Islas Turcas y Caicos',555,17,220,331,0,0,4),
    @('Tanzania',509,0,183,305,0,0,21),
    @('San Martin (Parte Holandesa)',495,13,302,174,0,0,19),
    @('Taiwan',489,0,471,11,0,0,7),
    @('Papua Nueva Guinea',471,0,232,234,0,0,5),
    @('Burundi',448,0,357,90,0,0,1),
    @('Comoras',448,21,410,31,0,0,7),
    @('Islas Feroe',412,1,379,33,0,0,0),
    @('Mauricio',356,0,335,11,0,0,10),
    @('Isla de Man',336,0,312,0,0,0,24),
    @('Eritrea',330,11,295,35,0,0,0),
    @('Mongolia',306,0,296,10,0,0,0),
    @('Gibraltar',298,3,246,52,0,0,0),
    @('Camboya',274,0,271,3,0,0,0),
    @('San Martin (Parte Francesa)',239,26,79,154,0,1,6),
    @('Butan',227,0,144,83,0,0,0),
    @('Islas Caimanes',205,0,204,0,0,0,1),
    @('Barbados',177,1,150,20,0,0,7),
    @('Bermudas',174,2,157,8,0,0,9),
    @('Brunei',145,1,139,3,0,0,3),
    @('Monaco',143,1,94,48,0,0,1),
    @('Seychelles',136,0,127,9,0,0,0),
    @('Liechtenstein',107,0,97,9,0,0,1),
    @('Antigua y Barbuda',95,1,91,1,0,0,3),
    @('Curazao',75,4,39,35,0,0,1),
    @('San Vicente y las Granadinas',61,0,57,4,0,0,0),
    @('Islas Virgenes Britanicas',47,0,8,38,0,0,1),
    @('Macao',46,0,46,0,0,0,0),
    @('Puerto Rico',39,0,1,36,0,0,2),
    @('Guam',32,0,0,31,0,0,1),
    @('Fiyi',29,0,23,4,0,0,2),
    @('Timor Oriental',27,0,25,2,0,0,0),
    @('Nueva Caledonia',26,3,23,3,0,0,0),
    @('Santa Lucia',26,0,26,0,0,0,0),
    @('Granada',24,0,24,0,0,0,0),
    @('Laos',22,0,21,1,0,0,0),
    @('Dominica',20,0,18,2,0,0,0),
    @('Bonaire, San Eustaquio y Saba',18,2,7,11,0,0,0),
    @('San Bartolome',18,0,13,5,0,0,0),
    @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
    @('San Cristobal y Nieves',17,0,17,0,0,0,0),
    @('Groenlandia',14,0,14,0,0,0,0),
    @('Islas Malvinas',13,0,13,0,0,0,0),
    @('Montserrat',13,0,12,0,0,0,1),
    @('Santa Sede',12,0,12,0,0,0,0),
    @('Sahara Occidental',10,0,8,1,0,0,1),
    @('San Pedro y Miquelon',5,0,4,1,0,0,0),
    @('Anguila',3,0,3,0,0,0,0)
)

$rows = $data.Count
$cols = 8
$arr = New-Object 'object[,]' $rows,$cols
for ($r = 0; $r -lt $rows; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $cols; $c++) {
        $arr[$r,$c] = $row[$c]
    }
}

$startRow = 4
$endRow = $startRow + $rows - 1
$ws.Range("A$startRow`:H$endRow").Value = $arr

"Rows written: " + $rows
"B4 (Estados Unidos): " + $ws.Range("B4").Value2
"H219 (Anguila): " + $ws.Range("H219").Value2
